$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A13's timestamp value (corrects fractional-day precision)
$ws.Cells.Item(13, 1).Value = 45875.66686240741

# Append new row 14 with the latest reading
$ws.Cells.Item(14, 1).Value = 45875.70852886758
$ws.Cells.Item(14, 2).Value = 2025
$ws.Cells.Item(14, 3).Value = 23
$ws.Cells.Item(14, 4).Value = 20.78
$ws.Cells.Item(14, 5).Value = 74.25
$ws.Cells.Item(14, 6).Value = 135.35
$ws.Cells.Item(14, 7).Value = 9.62
$ws.Cells.Item(14, 8).Value = "ESE"
$ws.Cells.Item(14, 9).Value = 0
$ws.Cells.Item(14, 10).Value = "17:00:16"

# Match the existing date-time number format / style used by column A
$ws.Cells.Item(14, 1).NumberFormat = $ws.Cells.Item(13, 1).NumberFormat
